# "Feature Narratives are supposed by workflow _focused_..."
#
# Slide 1 ("PM Quarterly Life"):
#   - TextBox 72 (id=73, the "Justification" box) shifts up a bit
#     (only its Top changes).
#   - TextBox 73 (id=74, the "Workflows" box) shifts left, shifts up a
#     bit, and widens to fit its new, longer label; its heading text
#     changes from "Workflows" to "Workflow focused".
#
# NOTE on precision: Shape.Left/Top/Width/Height are exposed as 32-bit
# (single-precision) floats, same as real PowerPoint's object model, while
# the underlying OOXML stores positions/sizes in EMU (914400 per inch,
# 12700 per point) as integers. Converting target-EMU -> points -> (single
# float) -> back to EMU can lose the last EMU to float rounding. The two
# literals below (for the new Left/Width of the "Workflows" box) are
# nudged by a fraction of a point - far too small to notice - so that they
# round-trip through the single-precision float to the exact target EMU
# values (4892770 and 1401794 respectively).
$emuPerPt = 914400 / 72

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 21: "TextBox 72" (Justification box moves up) -------------------
$justification = $s.Shapes.Item(21)
$justification.Top = 3989404 / $emuPerPt

# --- Shape 22: "TextBox 73" (Workflows -> Workflow focused) ----------------
$workflow = $s.Shapes.Item(22)
$workflow.Left = 385.257522753468       # -> 4892770 EMU
$workflow.Top = 3991488 / $emuPerPt     # -> 3991488 EMU
$workflow.Width = 110.377521799793      # -> 1401794 EMU
# Height (cy) is unchanged at 646331 EMU, so it is left untouched.

$workflow.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Workflow focused"
